$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2720.2
$ws.Range("I6").Value = 1250.5
$ws.Range("J6").Value = 3700
$ws.Range("K6").Value = 3751.5
$ws.Range("L6").Value = 11100
$ws.Range("M6").Value = -3639.5
$ws.Range("N6").Value = -11324
$ws.Range("H9").Value = 274.92307
$ws.Range("I9").Value = 121.28571
$ws.Range("J9").Value = 454.16666
$ws.Range("K9").Value = 121.28571
$ws.Range("L9").Value = 454.16666
$ws.Range("M9").Value = 47.71429000000001
$ws.Range("N9").Value = -792.16666
$ws.Range("H12").Value = 886.64
$ws.Range("I12").Value = 828.6
$ws.Range("K12").Value = 828.6
$ws.Range("M12").Value = -658.6
$ws.Range("H21").Value = 20539.75
$ws.Range("I21").Value = 11158.5
$ws.Range("J21").Value = 29921
$ws.Range("K21").Value = 11158.5
$ws.Range("L21").Value = 29921
$ws.Range("M21").Value = -10690.5
$ws.Range("N21").Value = -30857
$ws.Range("H23").Value = 20539.75
$ws.Range("I23").Value = 11158.5
$ws.Range("J23").Value = 29921
$ws.Range("K23").Value = 11158.5
$ws.Range("L23").Value = 29921
$ws.Range("M23").Value = -10924.5
$ws.Range("N23").Value = -30389
$ws.Range("H29").Value = 2754.6
$ws.Range("I29").Value = 1443.25
$ws.Range("J29").Value = 8000
$ws.Range("K29").Value = 4329.75
$ws.Range("L29").Value = 24000
$ws.Range("M29").Value = -4048.75
$ws.Range("N29").Value = -24562
$ws.Range("H33").Value = 1185.5834
$ws.Range("I33").Value = 1272.1818
$ws.Range("K33").Value = 1272.1818
$ws.Range("M33").Value = -1043.1818
$ws.Range("H38").Value = 1609.1111
$ws.Range("I38").Value = 1605.0769
$ws.Range("J38").Value = 1619.6
$ws.Range("K38").Value = 4815.2307
$ws.Range("L38").Value = 4858.799999999999
$ws.Range("M38").Value = -4443.2307
$ws.Range("N38").Value = -5602.799999999999
$ws.Range("H40").Value = 3879.8572
$ws.Range("I40").Value = 4329.5
$ws.Range("J40").Value = 2755.75
$ws.Range("K40").Value = 4329.5
$ws.Range("L40").Value = 2755.75
$ws.Range("M40").Value = -4154.5
$ws.Range("N40").Value = -3105.75
$ws.Range("H41").Value = 388.84616
$ws.Range("I41").Value = 368.625
$ws.Range("K41").Value = 368.625
$ws.Range("M41").Value = 71.375
$ws.Range("H43").Value = 957.53845
$ws.Range("I43").Value = 951.4
$ws.Range("J43").Value = 978
$ws.Range("K43").Value = 951.4
$ws.Range("L43").Value = 978
$ws.Range("M43").Value = -882.4
$ws.Range("N43").Value = -1116
$ws.Range("H58").Value = 3058.25
$ws.Range("I58").Value = 344.5
$ws.Range("J58").Value = 11199.5
$ws.Range("K58").Value = 1033.5
$ws.Range("L58").Value = 33598.5
$ws.Range("M58").Value = -883.5
$ws.Range("N58").Value = -33898.5
$ws.Range("H69").Value = 19864.416
$ws.Range("I69").Value = 7478.6
$ws.Range("J69").Value = 28711.428
$ws.Range("K69").Value = 22435.8
$ws.Range("L69").Value = 86134.284
$ws.Range("M69").Value = -21561.8
$ws.Range("N69").Value = -87882.284
$ws.Range("H72").Value = 19864.416
$ws.Range("I72").Value = 7478.6
$ws.Range("J72").Value = 28711.428
$ws.Range("K72").Value = 67307.40000000001
$ws.Range("L72").Value = 258402.852
$ws.Range("M72").Value = -62939.40000000001
$ws.Range("N72").Value = -267138.852
$ws.Range("H76").Value = 1292.7778
$ws.Range("I76").Value = 1292.7778
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 1292.7778
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -977.7778000000001
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 1292.7778
$ws.Range("I79").Value = 1292.7778
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 1292.7778
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -200.7778000000001
$ws.Range("N79").ClearContents()
$ws.Range("H80").Value = 2259.2856
$ws.Range("I80").Value = 2104.8
$ws.Range("J80").Value = 2399.7273
$ws.Range("K80").Value = 6314.400000000001
$ws.Range("L80").Value = 7199.1819
$ws.Range("M80").Value = -5316.400000000001
$ws.Range("N80").Value = -9195.1819
$ws.Range("H83").Value = 2259.2856
$ws.Range("I83").Value = 2104.8
$ws.Range("J83").Value = 2399.7273
$ws.Range("K83").Value = 18943.2
$ws.Range("L83").Value = 21597.5457
$ws.Range("M83").Value = -13951.2
$ws.Range("N83").Value = -31581.5457
$ws.Range("H92").Value = 840.58826
$ws.Range("I92").Value = 368
$ws.Range("J92").Value = 1974.8
$ws.Range("K92").Value = 368
$ws.Range("L92").Value = 1974.8
$ws.Range("M92").Value = 880
$ws.Range("N92").Value = -4470.8
$ws.Range("H94").Value = 23813836
$ws.Range("I94").Value = 23813836
$ws.Range("K94").Value = 23813836
$ws.Range("M94").Value = -23813385
$ws.Range("H98").Value = 2078.6
$ws.Range("J98").Value = 1899.5
$ws.Range("L98").Value = 1899.5
$ws.Range("N98").Value = -4895.5
$ws.Range("H106").Value = 100004080
$ws.Range("I106").Value = 100004080
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 100004080
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -100003449
$ws.Range("N106").ClearContents()
$ws.Range("H107").Value = 41668520
$ws.Range("I107").Value = 55557616
$ws.Range("J107").Value = 1232.3334
$ws.Range("K107").Value = 55557616
$ws.Range("L107").Value = 1232.3334
$ws.Range("M107").Value = -55555696
$ws.Range("N107").Value = -5072.3334
$ws.Range("H111").Value = 982.61536
$ws.Range("I111").Value = 768.5714
$ws.Range("J111").Value = 1232.3334
$ws.Range("K111").Value = 2305.7142
$ws.Range("L111").Value = 3697.0002
$ws.Range("M111").Value = 761.2857999999997
$ws.Range("N111").Value = -9831.0002
$ws.Range("H116").Value = 4577.7827
$ws.Range("I116").Value = 3752.6
$ws.Range("J116").Value = 4807
$ws.Range("K116").Value = 3752.6
$ws.Range("L116").Value = 4807
$ws.Range("M116").Value = -310.5999999999999
$ws.Range("N116").Value = -11691
$ws.Range("H122").Value = 2078.6
$ws.Range("J122").Value = 1899.5
$ws.Range("L122").Value = 5698.5
$ws.Range("N122").Value = -10598.5
$ws.Range("H131").Value = 6959
$ws.Range("I131").Value = 1703
$ws.Range("K131").Value = 5109
$ws.Range("M131").Value = -69
$ws.Range("H132").Value = 4537.2
$ws.Range("I132").Value = 3343.6345
$ws.Range("J132").Value = 9311.462
$ws.Range("K132").Value = 10030.9035
$ws.Range("L132").Value = 27934.386
$ws.Range("M132").Value = -7500.9035
$ws.Range("N132").Value = -32994.386
$ws.Range("H133").Value = 123333
$ws.Range("J133").Value = 123333
$ws.Range("L133").Value = 123333
$ws.Range("N133").Value = -133453
$ws.Range("H136").Value = 134999.73
$ws.Range("J136").Value = 134999.73
$ws.Range("L136").Value = 134999.73
$ws.Range("N136").Value = -145199.73
$ws.Range("H138").Value = 2056.5957
$ws.Range("I138").Value = 1106.1852
$ws.Range("J138").Value = 3339.65
$ws.Range("K138").Value = 3318.5556
$ws.Range("L138").Value = 10018.95
$ws.Range("M138").Value = 1821.4444
$ws.Range("N138").Value = -20298.95
$ws.Range("H141").Value = 11538.333
$ws.Range("I141").Value = 15789.125
$ws.Range("J141").Value = 3036.75
$ws.Range("K141").Value = 47367.375
$ws.Range("L141").Value = 9110.25
$ws.Range("M141").Value = -42187.375
$ws.Range("N141").Value = -19470.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 565.2308
$ws.Range("I4").Value = 227.14285
$ws.Range("K4").Value = 227.14285
$ws.Range("M4").Value = -111.14285
$ws.Range("H32").Value = 5170.625
$ws.Range("I32").Value = 2498.0881
$ws.Range("K32").Value = 2498.0881
$ws.Range("M32").Value = -2211.0881
$ws.Range("H43").Value = 16993.75
$ws.Range("J43").Value = 16993.75
$ws.Range("L43").Value = 16993.75
$ws.Range("N43").Value = -17619.75
$ws.Range("H45").Value = 33560.5
$ws.Range("I45").Value = 130000
$ws.Range("J45").Value = 1414
$ws.Range("K45").Value = 130000
$ws.Range("L45").Value = 1414
$ws.Range("M45").Value = -129623
$ws.Range("N45").Value = -2168
$ws.Range("H74").Value = 1509.174
$ws.Range("I74").Value = 1350.5
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 1350.5
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -476.5
$ws.Range("N74").Value = -6748
$ws.Range("H77").Value = 1509.174
$ws.Range("I77").Value = 1350.5
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 6752.5
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -2384.5
$ws.Range("N77").Value = -33736
$ws.Range("H88").Value = 4899.25
$ws.Range("J88").Value = 4866
$ws.Range("L88").Value = 4866
$ws.Range("N88").Value = -5678
$ws.Range("H91").Value = 4899.25
$ws.Range("J91").Value = 4866
$ws.Range("L91").Value = 4866
$ws.Range("N91").Value = -7674
$ws.Range("H102").Value = 6279.0557
$ws.Range("I102").Value = 5014
$ws.Range("J102").Value = 7860.375
$ws.Range("K102").Value = 5014
$ws.Range("L102").Value = 7860.375
$ws.Range("M102").Value = -3392
$ws.Range("N102").Value = -11104.375
$ws.Range("H112").Value = 279500
$ws.Range("J112").Value = 279500
$ws.Range("L112").Value = 279500
$ws.Range("N112").Value = -282454
$ws.Range("H122").Value = 1728.6757
$ws.Range("I122").Value = 1708.9062
$ws.Range("K122").Value = 5126.7186
$ws.Range("M122").Value = -2676.7186
$ws.Range("H132").Value = 2831.7778
$ws.Range("I132").Value = 2654.3333
$ws.Range("K132").Value = 7962.999899999999
$ws.Range("M132").Value = -5432.999899999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2050.4167
$ws.Range("I20").Value = 1972.619
$ws.Range("K20").Value = 1972.619
$ws.Range("M20").Value = -1725.619
$ws.Range("H22").Value = 374.47058
$ws.Range("I22").Value = 393.5625
$ws.Range("J22").Value = 69
$ws.Range("K22").Value = 393.5625
$ws.Range("L22").Value = 69
$ws.Range("M22").Value = -220.5625
$ws.Range("N22").Value = -415
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H55").Value = 97249
$ws.Range("J55").Value = 97249
$ws.Range("L55").Value = 97249
$ws.Range("N55").Value = -97795
$ws.Range("H86").Value = 1883.3334
$ws.Range("I86").Value = 1850
$ws.Range("K86").Value = 1850
$ws.Range("M86").Value = -727
$ws.Range("H89").Value = 1883.3334
$ws.Range("I89").Value = 1850
$ws.Range("K89").Value = 9250
$ws.Range("M89").Value = -3634
$ws.Range("H99").Value = 3449
$ws.Range("I99").Value = 4712.5
$ws.Range("J99").Value = 2185.5
$ws.Range("K99").Value = 4712.5
$ws.Range("L99").Value = 2185.5
$ws.Range("M99").Value = -3214.5
$ws.Range("N99").Value = -5181.5
$ws.Range("H105").Value = 4113.593
$ws.Range("I105").Value = 2655.8235
$ws.Range("K105").Value = 2655.8235
$ws.Range("M105").Value = -908.8235
$ws.Range("H120").Value = 90219.78
$ws.Range("J120").Value = 90219.78
$ws.Range("L120").Value = 90219.78
$ws.Range("N120").Value = -99895.78
$ws.Range("H134").Value = 3208.9363
$ws.Range("I134").Value = 2444.0789
$ws.Range("J134").Value = 6438.3335
$ws.Range("K134").Value = 7332.236699999999
$ws.Range("L134").Value = 19315.0005
$ws.Range("M134").Value = -4797.236699999999
$ws.Range("N134").Value = -24385.0005

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1074.5
$ws.Range("I16").Value = 1074.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1074.5
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -787.5
$ws.Range("N16").ClearContents()
$ws.Range("H22").Value = 935.2646999999999
$ws.Range("I22").Value = 695.6
$ws.Range("J22").Value = 1277.6428
$ws.Range("K22").Value = 695.6
$ws.Range("L22").Value = 1277.6428
$ws.Range("M22").Value = -345.6
$ws.Range("N22").Value = -1977.6428
$ws.Range("H31").Value = 1697.6666
$ws.Range("I31").Value = 1401.027
$ws.Range("J31").Value = 3069.625
$ws.Range("K31").Value = 1401.027
$ws.Range("L31").Value = 3069.625
$ws.Range("M31").Value = -1106.027
$ws.Range("N31").Value = -3659.625
$ws.Range("H34").Value = 1697.6666
$ws.Range("I34").Value = 1401.027
$ws.Range("J34").Value = 3069.625
$ws.Range("K34").Value = 1401.027
$ws.Range("L34").Value = 3069.625
$ws.Range("M34").Value = -1199.027
$ws.Range("N34").Value = -3473.625
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H48").Value = 74990
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 74990
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 74990
$ws.Range("M48").ClearContents()
$ws.Range("N48").Value = -75942
$ws.Range("H105").Value = 1344.8276
$ws.Range("I105").Value = 1434.7727
$ws.Range("J105").Value = 1062.1428
$ws.Range("K105").Value = 1434.7727
$ws.Range("L105").Value = 1062.1428
$ws.Range("M105").Value = 312.2273
$ws.Range("N105").Value = -4556.1428
$ws.Range("H113").Value = 1074.5
$ws.Range("I113").Value = 1074.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1074.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1095.5
$ws.Range("N113").ClearContents()
$ws.Range("H116").Value = 79999
$ws.Range("J116").Value = 79999
$ws.Range("L116").Value = 79999
$ws.Range("N116").Value = -89177
$ws.Range("H122").Value = 2741.5
$ws.Range("I122").Value = 2741.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8224.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5774.5
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 5873.6177
$ws.Range("I132").Value = 6150.1665
$ws.Range("J132").Value = 3799.5
$ws.Range("K132").Value = 18450.4995
$ws.Range("L132").Value = 11398.5
$ws.Range("M132").Value = -15920.4995
$ws.Range("N132").Value = -16458.5
$ws.Range("H134").Value = 6083.1763
$ws.Range("I134").Value = 5147.375
$ws.Range("K134").Value = 15442.125
$ws.Range("M134").Value = -12907.125

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 83.63636
$ws.Range("I2").Value = 30.5
$ws.Range("J2").Value = 114
$ws.Range("K2").Value = 183
$ws.Range("L2").Value = 684
$ws.Range("M2").Value = -70
$ws.Range("N2").Value = -910
$ws.Range("H38").Value = 115
$ws.Range("I38").Value = 89.38461
$ws.Range("J38").Value = 181.6
$ws.Range("K38").Value = 268.15383
$ws.Range("L38").Value = 544.8
$ws.Range("M38").Value = 78.84617000000003
$ws.Range("N38").Value = -1238.8
$ws.Range("H46").Value = 430.25
$ws.Range("I46").Value = 273.75
$ws.Range("K46").Value = 821.25
$ws.Range("M46").Value = -730.25
$ws.Range("H50").Value = 3543.0908
$ws.Range("I50").Value = 2498.5
$ws.Range("J50").Value = 3775.2222
$ws.Range("K50").Value = 7495.5
$ws.Range("L50").Value = 11325.6666
$ws.Range("M50").Value = -7014.5
$ws.Range("N50").Value = -12287.6666
$ws.Range("H53").Value = 3543.0908
$ws.Range("I53").Value = 2498.5
$ws.Range("J53").Value = 3775.2222
$ws.Range("K53").Value = 7495.5
$ws.Range("L53").Value = 11325.6666
$ws.Range("M53").Value = -7014.5
$ws.Range("N53").Value = -12287.6666
$ws.Range("H68").Value = 125014130
$ws.Range("I68").Value = 1002
$ws.Range("J68").Value = 166685170
$ws.Range("K68").Value = 3006
$ws.Range("L68").Value = 500055510
$ws.Range("M68").Value = -2195
$ws.Range("N68").Value = -500057132
$ws.Range("H71").Value = 125014130
$ws.Range("I71").Value = 1002
$ws.Range("J71").Value = 166685170
$ws.Range("K71").Value = 9018
$ws.Range("L71").Value = 1500166530
$ws.Range("M71").Value = -4962
$ws.Range("N71").Value = -1500174642
$ws.Range("H99").Value = 2518.8
$ws.Range("I99").Value = 1183.1428
$ws.Range("K99").Value = 3549.4284
$ws.Range("M99").Value = -1303.4284
$ws.Range("H121").Value = 15153118
$ws.Range("I121").Value = 33333510
$ws.Range("J121").Value = 2791.6667
$ws.Range("K121").Value = 100000530
$ws.Range("L121").Value = 8375.000100000001
$ws.Range("M121").Value = -99999220
$ws.Range("N121").Value = -10995.0001
$ws.Range("H132").Value = 1770.1578
$ws.Range("I132").Value = 1726.3
$ws.Range("J132").Value = 1818.8889
$ws.Range("K132").Value = 15536.7
$ws.Range("L132").Value = 16370.0001
$ws.Range("M132").Value = -13006.7
$ws.Range("N132").Value = -21430.0001
$ws.Range("H133").Value = 2242.3572
$ws.Range("I133").Value = 2242.3572
$ws.Range("K133").Value = 6727.071599999999
$ws.Range("M133").Value = -1667.071599999999
$ws.Range("H138").Value = 2777.4285
$ws.Range("I138").Value = 2777.4285
$ws.Range("K138").Value = 8332.2855
$ws.Range("M138").Value = -3192.2855
$ws.Range("H140").Value = 1250.75
$ws.Range("I140").Value = 990
$ws.Range("K140").Value = 2970
$ws.Range("M140").Value = 2210

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 111111360
$ws.Range("I2").Value = 155555740
$ws.Range("J2").Value = 372.5
$ws.Range("K2").Value = 155555740
$ws.Range("L2").Value = 372.5
$ws.Range("M2").Value = -155555627
$ws.Range("N2").Value = -598.5
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H70").Value = 16057.667
$ws.Range("I70").Value = 19848
$ws.Range("K70").Value = 19848
$ws.Range("M70").Value = -19578
$ws.Range("H73").Value = 16057.667
$ws.Range("I73").Value = 19848
$ws.Range("K73").Value = 19848
$ws.Range("M73").Value = -18912
$ws.Range("H80").Value = 2648.761
$ws.Range("J80").Value = 2110.158
$ws.Range("L80").Value = 2110.158
$ws.Range("N80").Value = -4106.157999999999
$ws.Range("H83").Value = 2648.761
$ws.Range("J83").Value = 2110.158
$ws.Range("L83").Value = 10550.79
$ws.Range("N83").Value = -20534.79
$ws.Range("H97").Value = 2258.2856
$ws.Range("I97").Value = 1105.5
$ws.Range("J97").Value = 3795.3333
$ws.Range("K97").Value = 1105.5
$ws.Range("L97").Value = 3795.3333
$ws.Range("M97").Value = -609.5
$ws.Range("N97").Value = -4787.3333
$ws.Range("H102").Value = 5309.6665
$ws.Range("I102").Value = 6104.6665
$ws.Range("J102").Value = 4249.6665
$ws.Range("K102").Value = 6104.6665
$ws.Range("L102").Value = 4249.6665
$ws.Range("M102").Value = -4482.6665
$ws.Range("N102").Value = -7493.6665
$ws.Range("H122").Value = 479813.72
$ws.Range("I122").Value = 835947.7
$ws.Range("K122").Value = 2507843.1
$ws.Range("M122").Value = -2505393.1
$ws.Range("H132").Value = 5972.7
$ws.Range("I132").Value = 5480.5
$ws.Range("J132").Value = 6095.75
$ws.Range("K132").Value = 16441.5
$ws.Range("L132").Value = 18287.25
$ws.Range("M132").Value = -13911.5
$ws.Range("N132").Value = -23347.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3807.1614
$ws.Range("I7").Value = 3384.1304
$ws.Range("J7").Value = 5023.375
$ws.Range("K7").Value = 3384.1304
$ws.Range("L7").Value = 5023.375
$ws.Range("M7").Value = -3272.1304
$ws.Range("N7").Value = -5247.375
$ws.Range("H16").Value = 2548
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H46").Value = 2884.1155
$ws.Range("J46").Value = 3294.6191
$ws.Range("L46").Value = 3294.6191
$ws.Range("N46").Value = -3670.6191
$ws.Range("H55").Value = 1092.72
$ws.Range("J55").Value = 1404.5714
$ws.Range("L55").Value = 1404.5714
$ws.Range("N55").Value = -1750.5714
$ws.Range("H93").Value = 1775
$ws.Range("I93").Value = 1836.5
$ws.Range("J93").Value = 1375.25
$ws.Range("K93").Value = 1836.5
$ws.Range("L93").Value = 1375.25
$ws.Range("M93").Value = -588.5
$ws.Range("N93").Value = -3871.25
$ws.Range("H100").Value = 3000
$ws.Range("I100").Value = 3000
$ws.Range("K100").Value = 3000
$ws.Range("M100").Value = -2459
$ws.Range("H122").Value = 5744.3125
$ws.Range("I122").Value = 3701.3
$ws.Range("K122").Value = 11103.9
$ws.Range("M122").Value = -8653.900000000001
$ws.Range("H126").Value = 3807.1614
$ws.Range("I126").Value = 3384.1304
$ws.Range("J126").Value = 5023.375
$ws.Range("K126").Value = 10152.3912
$ws.Range("L126").Value = 15070.125
$ws.Range("M126").Value = -7682.3912
$ws.Range("N126").Value = -20010.125
$ws.Range("H132").Value = 45724.16
$ws.Range("I132").Value = 53709.715
$ws.Range("J132").Value = 3800
$ws.Range("K132").Value = 161129.145
$ws.Range("L132").Value = 11400
$ws.Range("M132").Value = -158599.145
$ws.Range("N132").Value = -16460
$ws.Range("H136").Value = 6440005.5
$ws.Range("I136").Value = 8578296
$ws.Range("J136").Value = 25132.715
$ws.Range("K136").Value = 25734888
$ws.Range("L136").Value = 75398.145
$ws.Range("M136").Value = -25732338
$ws.Range("N136").Value = -80498.145

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 23999.4
$ws.Range("I2").Value = 19999
$ws.Range("J2").Value = 26666.334
$ws.Range("K2").Value = 19999
$ws.Range("L2").Value = 26666.334
$ws.Range("M2").Value = -19887
$ws.Range("N2").Value = -26890.334
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H113").Value = 1650.6666
$ws.Range("I113").Value = 999.5
$ws.Range("J113").Value = 1976.25
$ws.Range("K113").Value = 2998.5
$ws.Range("L113").Value = 5928.75
$ws.Range("M113").Value = -828.5
$ws.Range("N113").Value = -10268.75
$ws.Range("H126").Value = 5149.5864
$ws.Range("I126").Value = 4541.4
$ws.Range("K126").Value = 13624.2
$ws.Range("M126").Value = -11154.2
$ws.Range("H132").Value = 2384.7778
$ws.Range("I132").Value = 2066.8572
$ws.Range("K132").Value = 6200.571599999999
$ws.Range("M132").Value = -3670.571599999999
